$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (issue number & week range) ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Fix up cell styles (number <-> text) for cells changing data type, by
#     copying format+value from a same-style donor cell before writing the
#     real target value. This keeps style indices aligned with the donor.
$ws.Range("G14").Copy($ws.Range("C15"))
$ws.Range("G14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("G14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("G26"))
$ws.Range("E14").Copy($ws.Range("H26"))
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))

# --- Write final cell values ---
$ws.Range("N14").Value = -50
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = "0"
$ws.Range("H15").Value = "***.*"
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 60
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -11.111111111111
$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = 29.268292682926
$ws.Range("L16").Value = 43.243243243243
$ws.Range("M16").Value = -17.1875
$ws.Range("N16").Value = -84.084084084084
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -47.058823529411
$ws.Range("I17").Value = 71
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 24.561403508771
$ws.Range("L17").Value = 5.970149253731
$ws.Range("M17").Value = 69.047619047619
$ws.Range("N17").Value = -36.036036036036
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -17.647058823529
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = 35.714285714285
$ws.Range("L18").Value = 63.793103448275
$ws.Range("M18").Value = -11.214953271028
$ws.Range("N18").Value = -84.140233722871
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -7.692307692307
$ws.Range("I19").Value = 239
$ws.Range("J19").Value = 275
$ws.Range("K19").Value = -13.090909090909
$ws.Range("L19").Value = 14.903846153846
$ws.Range("M19").Value = 57.236842105263
$ws.Range("N19").Value = -12.454212454212
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 53
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = 39.473684210526
$ws.Range("L20").Value = 140.909090909091
$ws.Range("M20").Value = -25.352112676056
$ws.Range("N20").Value = -92.274052478134
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 19.230769230769
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -12.14953271028
$ws.Range("I21").Value = 520
$ws.Range("J21").Value = 493
$ws.Range("K21").Value = 5.476673427991
$ws.Range("L21").Value = 30.653266331658
$ws.Range("M21").Value = 16.853932584269
$ws.Range("N21").Value = -74.167908594138
$ws.Range("M22").Value = -76.923076923076
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 12
$ws.Range("F24").Value = 147
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 18.548387096774
$ws.Range("I24").Value = 681
$ws.Range("J24").Value = 613
$ws.Range("K24").Value = 11.092985318107
$ws.Range("L24").Value = 63.309352517985
$ws.Range("M24").Value = 87.087912087912
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 150
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 103.571428571429
$ws.Range("I25").Value = 224
$ws.Range("J25").Value = 157
$ws.Range("K25").Value = 42.675159235668
$ws.Range("L25").Value = 34.939759036144
$ws.Range("M25").Value = 25.139664804469
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = "0"
$ws.Range("H26").Value = "***.*"
$ws.Range("I26").Value = 10
$ws.Range("K26").Value = -9.090909090909
$ws.Range("L26").Value = 100
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 80
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -10
$ws.Range("L27").Value = 5.882352941176
$ws.Range("G30").Value = "0"
$ws.Range("H30").Value = "***.*"
